$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose content would otherwise be auto-converted
# by Excel (numeric-looking quantities, and date-like strings) so they are
# stored as text, matching the source data.
$ws.Range("I7:I11").NumberFormat = "@"
$ws.Range("Y7:Y11").NumberFormat = "@"
$ws.Range("AA7:AA11").NumberFormat = "@"

# Row 7
$ws.Range("A7").Value = 131106436
$ws.Range("B7").Value = 5493
$ws.Range("D7").Value = 'NT'
$ws.Range("E7").Value = 101410
$ws.Range("F7").Value = 'Reliktbock'
$ws.Range("G7").Value = 'Nothorhina muricata'
$ws.Range("H7").Value = '(Dalman, 1817)'
$ws.Range("I7").Value = '2'
$ws.Range("J7").Value = 'ex.'
$ws.Range("P7").Value = 'Svartmyran, Mpd'
$ws.Range("Q7").Value = 616762
$ws.Range("R7").Value = 6934714
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = 'Västernorrland'
$ws.Range("U7").Value = 'Timrå'
$ws.Range("V7").Value = 'Medelpad'
$ws.Range("W7").Value = 'Timrå'
$ws.Range("X7").Value = '2025_0743'
$ws.Range("Y7").Value = '2025-07-02'
$ws.Range("Z7").Value = '11:39'
$ws.Range("AA7").Value = '2025-07-02'
$ws.Range("AB7").Value = '11:39'
$ws.Range("AC7").Value = 'Två kläckhål'
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AT7").Value = "'"
$ws.Range("AW7").Value = 'David Isaksson'
$ws.Range("AX7").Value = 'David Isaksson'
$ws.Range("AY7").Value = 'Kustpaketet'

# Row 8
$ws.Range("A8").Value = 131108352
$ws.Range("B8").Value = 80214
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 388
$ws.Range("F8").Value = 'Stiftgelélav'
$ws.Range("G8").Value = 'Collema furfuraceum'
$ws.Range("H8").Value = '(Arnold) Du Rietz'
$ws.Range("I8").Value = '1'
$ws.Range("J8").Value = 'bålar'
$ws.Range("P8").Value = 'S Svartmyran, Mpd'
$ws.Range("Q8").Value = 616863
$ws.Range("R8").Value = 6934788
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = 'Västernorrland'
$ws.Range("U8").Value = 'Timrå'
$ws.Range("V8").Value = 'Medelpad'
$ws.Range("W8").Value = 'Timrå'
$ws.Range("X8").Value = '2025_0758'
$ws.Range("Y8").Value = '2025-07-02'
$ws.Range("Z8").Value = '14:47'
$ws.Range("AA8").Value = '2025-07-02'
$ws.Range("AB8").Value = '14:47'
$ws.Range("AC8").Value = 'Asp'
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AT8").Value = "'"
$ws.Range("AW8").Value = 'David Isaksson'
$ws.Range("AX8").Value = 'Måns Svensson'
$ws.Range("AY8").Value = 'Kustpaketet'

# Row 9
$ws.Range("A9").Value = 131106422
$ws.Range("B9").Value = 80252
$ws.Range("D9").Value = 'LC'
$ws.Range("E9").Value = 6456
$ws.Range("F9").Value = 'Skinnlav'
$ws.Range("G9").Value = 'Leptogium saturninum'
$ws.Range("H9").Value = '(Dicks.) Nyl.'
$ws.Range("I9").Value = '3'
$ws.Range("J9").Value = 'bålar'
$ws.Range("P9").Value = 'S Svartmyran, Mpd'
$ws.Range("Q9").Value = 616854
$ws.Range("R9").Value = 6934782
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = 'Västernorrland'
$ws.Range("U9").Value = 'Timrå'
$ws.Range("V9").Value = 'Medelpad'
$ws.Range("W9").Value = 'Timrå'
$ws.Range("X9").Value = '2025_0757'
$ws.Range("Y9").Value = '2025-07-02'
$ws.Range("Z9").Value = '14:43'
$ws.Range("AA9").Value = '2025-07-02'
$ws.Range("AB9").Value = '14:43'
$ws.Range("AC9").Value = 'Asp'
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AT9").Value = "'"
$ws.Range("AW9").Value = 'David Isaksson'
$ws.Range("AX9").Value = 'Måns Svensson'
$ws.Range("AY9").Value = 'Kustpaketet'

# Row 10
$ws.Range("A10").Value = 131106420
$ws.Range("B10").Value = 80377
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 6462
$ws.Range("F10").Value = 'Stuplav'
$ws.Range("G10").Value = 'Nephroma bellum'
$ws.Range("H10").Value = '(Spreng.) Tuck.'
$ws.Range("I10").Value = '1'
$ws.Range("J10").Value = 'dm²'
$ws.Range("P10").Value = 'S Svartmyran, Mpd'
$ws.Range("Q10").Value = 616876
$ws.Range("R10").Value = 6934813
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 'Västernorrland'
$ws.Range("U10").Value = 'Timrå'
$ws.Range("V10").Value = 'Medelpad'
$ws.Range("W10").Value = 'Timrå'
$ws.Range("X10").Value = '2025_0759'
$ws.Range("Y10").Value = '2025-07-02'
$ws.Range("Z10").Value = '14:51'
$ws.Range("AA10").Value = '2025-07-02'
$ws.Range("AB10").Value = '14:51'
$ws.Range("AC10").Value = 'Asp'
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AT10").Value = "'"
$ws.Range("AW10").Value = 'David Isaksson'
$ws.Range("AX10").Value = 'Måns Svensson'
$ws.Range("AY10").Value = 'Kustpaketet'

# Row 11
$ws.Range("A11").Value = 131106423
$ws.Range("B11").Value = 80348
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 6458
$ws.Range("F11").Value = 'Lunglav'
$ws.Range("G11").Value = 'Lobaria pulmonaria'
$ws.Range("H11").Value = '(L.) Hoffm.'
$ws.Range("I11").Value = '1'
$ws.Range("J11").Value = 'bålar'
$ws.Range("P11").Value = 'S Svartmyran, Mpd'
$ws.Range("Q11").Value = 616860
$ws.Range("R11").Value = 6934783
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = 'Västernorrland'
$ws.Range("U11").Value = 'Timrå'
$ws.Range("V11").Value = 'Medelpad'
$ws.Range("W11").Value = 'Timrå'
$ws.Range("X11").Value = '2025_0756'
$ws.Range("Y11").Value = '2025-07-02'
$ws.Range("Z11").Value = '14:43'
$ws.Range("AA11").Value = '2025-07-02'
$ws.Range("AB11").Value = '14:43'
$ws.Range("AC11").Value = 'På asp'
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AT11").Value = "'"
$ws.Range("AW11").Value = 'David Isaksson'
$ws.Range("AX11").Value = 'Måns Svensson'
$ws.Range("AY11").Value = 'Kustpaketet'
